$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)

foreach ($row in 2..4) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # Column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
